$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" standalone="yes"?>' + `
             '<?mso-application progid="Word.Document"?>' + `
             '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
             '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
             '<pkg:xmlData>' + `
             '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
             '<w:body><w:p>'
$xmlFooter = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Apply-Para($Index, $Inner) {
    $p = $d.Paragraphs.Item($Index)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $before = $r.Text
    $r.InsertXML($xmlHeader + $Inner + $xmlFooter)
    Write-Output ("Para {0}: [{1}] -> [{2}]" -f $Index, $before, $r.Text)
}

Apply-Para 4 '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Set up and initialize chips of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Espressif</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Apply-Para 5 '<w:r><w:t xml:space="preserve">Use “pip install </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esptool</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” to get installation tool for chips produced by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Espressif</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (esp8266, esp32, …)</w:t></w:r>'
Apply-Para 6 '<w:r><w:t xml:space="preserve">Use “esptool.py &lt;arguments…&gt;” in CMD to install </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>micropython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> firmware for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> devices mentioned above.</w:t></w:r>'
Apply-Para 9 '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>I</w:t></w:r><w:r><w:t xml:space="preserve">f </w:t></w:r><w:r><w:t>encounter</w:t></w:r><w:r><w:t xml:space="preserve"> “failed to create process” when try to use esptool.py, download </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esptool</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> project files on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and directly run esptool.py in project directory through command “python </w:t></w:r><w:r><w:t>&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>absolutepath</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;</w:t></w:r><w:r><w:t>\...\</w:t></w:r><w:r><w:t>esptool.py”.</w:t></w:r>'
Apply-Para 11 '<w:r><w:t xml:space="preserve">Use “import </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” to get hardware information of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> devices in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>micropython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> shell.</w:t></w:r>'
Apply-Para 12 '<w:r><w:t xml:space="preserve">E.g. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esp.flash_size</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() returns the size of on-board flash </w:t></w:r><w:r><w:t>in</w:t></w:r><w:r><w:t xml:space="preserve"> byte</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t>.</w:t></w:r>'
Apply-Para 13 '<w:r><w:t xml:space="preserve">Use help(obj) in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>micropython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> shell to get information of functions, classes, variables and other objects in modules.</w:t></w:r>'
Apply-Para 14 '<w:r><w:t>E.g. help(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>machine.SPI</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">) returns names of methods and properties of SPI class in module machine. </w:t></w:r>'
Apply-Para 48 '<w:proofErr w:type="spellStart"/><w:r><w:t>Guidence</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>M</w:t></w:r><w:r><w:t>icropython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> SPI class</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>:</w:t></w:r>'
Apply-Para 54 '<w:r><w:t xml:space="preserve">So before writing data, GPIO connected to LOAD should </w:t></w:r><w:r><w:t>sit at</w:t></w:r><w:r><w:t xml:space="preserve"> low.</w:t></w:r>'
Apply-Para 62 '<w:r><w:t>H</w:t></w:r><w:r><w:t xml:space="preserve">ence, </w:t></w:r><w:r><w:t xml:space="preserve">for this chip, </w:t></w:r><w:r><w:t>polarity and phase argument</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> of SPI class shou</w:t></w:r><w:r><w:t>ld</w:t></w:r><w:r><w:t xml:space="preserve"> be 0.</w:t></w:r>'
Apply-Para 63 '<w:r><w:t>Before writing data for LEDs, default values need to be submitted to other registers so that max7219 can run stably</w:t></w:r><w:r><w:t xml:space="preserve">, including decode mode, intensity, scan limit, operation mode, display test. More detailed </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">information </w:t></w:r><w:r><w:t>is</w:t></w:r><w:r><w:t xml:space="preserve"> presented in the datasheet.</w:t></w:r>'
Apply-Para 65 '<w:r><w:t>Priority of effect of display test is higher than that of shutdown</w:t></w:r>'

Write-Output "ALL DONE"
